$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content so the shared-strings table and grid are rebuilt cleanly
$ws.Cells.Clear() | Out-Null

# Column widths: A=34, C:E=21.83203125, F=10.83203125
# (values below are pre-compensated for this engine's character->pixel
# rounding in ColumnWidth so the stored XML "width" lands on target)
$ws.Columns("A").ColumnWidth = 33.285714285714285
$ws.Columns("C:E").ColumnWidth = 21.142857142857142
$ws.Columns("F").ColumnWidth = 10.142857142857142

# Populate cells per final layout
$ws.Range("A1").Value = "Route"
$ws.Range("B1").Value = "Method"
$ws.Range("C1").Value = "Params"
$ws.Range("F1").Value = "Response"
$ws.Range("A2").Value = "/signup"
$ws.Range("B2").Value = "POST"
$ws.Range("C2").Value = "user_name"
$ws.Range("F2").Value = "id"
$ws.Range("C3").Value = "password"
$ws.Range("A5").Value = "/login"
$ws.Range("B5").Value = "POST"
$ws.Range("C5").Value = "user_name"
$ws.Range("F5").Value = "id"
$ws.Range("C6").Value = "password"
$ws.Range("A8").Value = "/events"
$ws.Range("B8").Value = "GET"
$ws.Range("F8").Value = "events"
$ws.Range("G8").Value = "[{"
$ws.Range("G9").Value = "eventID"
$ws.Range("G10").Value = "type"
$ws.Range("G11").Value = "user"
$ws.Range("H11").Value = "{"
$ws.Range("H12").Value = "username"
$ws.Range("H13").Value = "profilePhotoURL"
$ws.Range("H14").Value = "}"
$ws.Range("G15").Value = "item"
$ws.Range("H15").Value = "{"
$ws.Range("H16").Value = "photoURL"
$ws.Range("H17").Value = "}"
$ws.Range("G18").Value = "poll"
$ws.Range("H18").Value = "{"
$ws.Range("H19").Value = "pollID"
$ws.Range("H20").Value = "title"
$ws.Range("H21").Value = "owner"
$ws.Range("I21").Value = "{"
$ws.Range("I22").Value = "username"
$ws.Range("I23").Value = "}"
$ws.Range("H24").Value = "}"
$ws.Range("G25").Value = "}"
$ws.Range("G26").Value = ",…"
$ws.Range("G27").Value = "]"
$ws.Range("A28").Value = "/events"
$ws.Range("B28").Value = "POST"
$ws.Range("C28").Value = "type"
$ws.Range("C29").Value = "user "
$ws.Range("D29").Value = "{ id }"
$ws.Range("C30").Value = "poll"
$ws.Range("D30").Value = "{ pollID }"
$ws.Range("C31").Value = "item"
$ws.Range("D31").Value = "{ }"
$ws.Range("A33").Value = "/polls"
$ws.Range("B33").Value = "POST"
$ws.Range("C33").Value = "title"
$ws.Range("F33").Value = "pollID"
$ws.Range("C34").Value = "ownerID"
$ws.Range("F34").Value = "ownerID"
$ws.Range("C35").Value = "state (default = `"EDITING`")"
$ws.Range("C36").Value = "totalVotes (default = 0)"
$ws.Range("C37").Value = "maxVotesForSingleItem (default = 1)"
$ws.Range("C38").Value = "startTime"
$ws.Range("A40").Value = "/polls/:pollID"
$ws.Range("B40").Value = "GET"
$ws.Range("C40").Value = "pollID"

# Freeze header row and set active selection
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F6").Select() | Out-Null
